# Fri Sep  6 19:52:02 UTC 2024 cryptos refresh: update Price/Volume(1h)
# columns for every coin row, and swap the USDe / PancakeSwap rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'53.591.18"
$ws.Range("E2").Value = "  -4.36%  "

$ws.Range("D3").Value = "'2.216.49"
$ws.Range("E3").Value = "  -5.96%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'484.08"
$ws.Range("E5").Value = "  -3.12%  "

$ws.Range("D6").Value = "'125.30"
$ws.Range("E6").Value = "  -2.70%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  -5.01%  "

$ws.Range("D9").Value = "'2.210.69"
$ws.Range("E9").Value = "  -6.14%  "

$ws.Range("D10").Value = "'0.0909"
$ws.Range("E10").Value = "  -6.67%  "

$ws.Range("E11").Value = "  -1.47%  "

$ws.Range("D12").Value = "'4.64"
$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").Value = "'0.313"
$ws.Range("E13").Value = "  -2.72%  "

$ws.Range("D14").Value = "'2.609.13"
$ws.Range("E14").Value = "  -5.94%  "

$ws.Range("D15").Value = "'20.90"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").Value = "'53.523.55"
$ws.Range("E16").Value = "  -4.42%  "

$ws.Range("E17").Value = "  -3.25%  "

$ws.Range("D18").Value = "'2.215.20"
$ws.Range("E18").Value = "  -15.37%  "

$ws.Range("D19").Value = "'9.51"
$ws.Range("E19").Value = "  -4.47%  "

$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").Value = "'297.17"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("D22").Value = "'6.08"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'63.10"
$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("E27").Value = "  -2.50%  "

$ws.Range("D28").Value = "'6.94"
$ws.Range("E28").Value = "  -3.28%  "

$ws.Range("D29").Value = "'169.14"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.57"
$ws.Range("E30").Value = "  -3.52%  "

$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").Value = "'0.0₃0671"
$ws.Range("E32").Value = "  -4.98%  "

$ws.Range("D33").Value = "'0.993"
$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("D35").Value = "'1.05"
$ws.Range("E35").Value = "  -2.97%  "

$ws.Range("D36").Value = "'17.32"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("E37").Value = "  -2.07%  "

$ws.Range("D38").Value = "'0.825"
$ws.Range("E38").Value = "  +5.32%  "

$ws.Range("D39").Value = "'3.55"
$ws.Range("E39").Value = "  -4.74%  "

$ws.Range("D40").Value = "'35.68"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("E43").Value = "  -2.60%  "

$ws.Range("D44").Value = "'122.39"
$ws.Range("E44").Value = "  -5.12%  "

$ws.Range("D45").Value = "'4.62"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").Value = "'0.0876"
$ws.Range("E46").Value = "  -2.66%  "

$ws.Range("D47").Value = "'0.525"
$ws.Range("E47").Value = "  -6.06%  "

$ws.Range("D48").Value = "'228.52"
$ws.Range("E48").Value = "  -4.10%  "

$ws.Range("D49").Value = "'0.0467"
$ws.Range("E49").Value = "  -2.31%  "

$ws.Range("D50").Value = "'0.0200"
$ws.Range("E50").Value = "  -2.75%  "

$ws.Range("D51").Value = "'15.89"
$ws.Range("E51").Value = "  -4.76%  "
